$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/value updates (safe -- Excel will not misparse these as numbers)
$plainUpdates = @{
    "D2" = "43.054.18"
    "E2" = "  +0.14%  "
    "D3" = "2.300.61"
    "E3" = "  +0.28%  "
    "E4" = "  +0.03%  "
    "E5" = "  -0.10%  "
    "E6" = "  -1.35%  "
    "E7" = "  +3.65%  "
    "E8" = "  -0.01%  "
    "E9" = "  +1.45%  "
    "E10" = "  +0.30%  "
    "E11" = "  +0.69%  "
    "E13" = "  -1.88%  "
    "E14" = "  -0.54%  "
    "D15" = "2.657.91"
    "E15" = "  +0.23%  "
    "D16" = "2.307.30"
    "E16" = "  +1.44%  "
    "E17" = "  -1.17%  "
    "D18" = "42.939.21"
    "E18" = "  +0.13%  "
    "E19" = "  +4.07%  "
    "D20" = "0.0₃0913"
    "E20" = "  +1.36%  "
    "E21" = "  +0.60%  "
    "E22" = "  +0.84%  "
    "E23" = "  +1.00%  "
    "E24" = "  -0.91%  "
    "E25" = "  -0.45%  "
    "E26" = "  -0.63%  "
    "E27" = "  -0.11%  "
    "E29" = "  -12.68%  "
    "E30" = "  +0.53%  "
    "E31" = "  -2.07%  "
    "E32" = "  -4.12%  "
    "E33" = "  +0.04%  "
    "E34" = "  +2.24%  "
    "E35" = "  +2.84%  "
    "E36" = "  +2.64%  "
    "E37" = "  +0.40%  "
    "E38" = "  +1.70%  "
    "E39" = "  +0.55%  "
    "E40" = "  -0.19%  "
    "E41" = "  +1.90%  "
    "E42" = "  -1.42%  "
    "D43" = "2.015.43"
    "E43" = "  +2.59%  "
    "E44" = "  -1.55%  "
    "E45" = "  -3.96%  "
    "E46" = "  +1.49%  "
    "E47" = "  +0.37%  "
    "E48" = "  -1.47%  "
    "E49" = "  -1.85%  "
    "D50" = "2.530.72"
    "E50" = "  +0.45%  "
    "E51" = "  -0.57%  "
}

foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}

# Numeric-looking price strings: must stay text (column D uses "." as a
# thousands separator, e.g. "43.054.18"), so force text entry the same way
# Excel does for a user-typed leading apostrophe, then strip the resulting
# quote-prefix cell style back off so formatting is untouched.
$forcedTextUpdates = @{
    "D5" = "299.97"
    "D6" = "97.79"
    "D10" = "36.30"
    "D13" = "17.80"
    "D17" = "0.789"
    "D22" = "68.27"
    "D23" = "237.93"
    "D31" = "163.49"
    "D32" = "33.04"
    "D33" = "1.00"
    "D35" = "18.16"
    "D36" = "4.74"
    "D45" = "2.22"
    "D47" = "17.48"
    "D48" = "2.84"
}

foreach ($addr in $forcedTextUpdates.Keys) {
    $ws.Range($addr).Value = "'" + $forcedTextUpdates[$addr]
    $ws.Range($addr).ClearFormats()
}
